$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ntn1"
$ws.Cells.Item(2,3).Value = "Adora2b"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.8750386666666667
$ws.Cells.Item(2,8).Value = 2.625116
$ws.Cells.Item(2,9).Value = 0.05304058862308838
$ws.Cells.Item(2,10).Value = 0.05304058862308838
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 3.888768333333333
$ws.Cells.Item(2,14).Value = 11.666305
$ws.Cells.Item(2,15).Value = 0.2727944290034114
$ws.Cells.Item(2,16).Value = 0.2727944290034114
$ws.Cells.Item(2,17).Value = 3.402822657375556
$ws.Cells.Item(2,18).Value = 30.62540391638
$ws.Cells.Item(2,19).Value = 0.01446917708744023
$ws.Cells.Item(2,20).Value = 0.01446917708744023

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ntn1"
$ws.Cells.Item(3,3).Value = "Adora2b"
$ws.Cells.Item(3,4).Value = "M2"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.8750386666666667
$ws.Cells.Item(3,8).Value = 2.625116
$ws.Cells.Item(3,9).Value = 0.05304058862308838
$ws.Cells.Item(3,10).Value = 0.05304058862308838
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 7.624472
$ws.Cells.Item(3,14).Value = 22.873416
$ws.Cells.Item(3,15).Value = 0.5348514767167063
$ws.Cells.Item(3,16).Value = 0.5348514767167063
$ws.Cells.Item(3,17).Value = 6.671707812917334
$ws.Cells.Item(3,18).Value = 60.045370316256
$ws.Cells.Item(3,19).Value = 0.02836883715098215
$ws.Cells.Item(3,20).Value = 0.02836883715098215

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ntn1"
$ws.Cells.Item(4,3).Value = "Adora2b"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.8750386666666667
$ws.Cells.Item(4,8).Value = 2.625116
$ws.Cells.Item(4,9).Value = 0.05304058862308838
$ws.Cells.Item(4,10).Value = 0.05304058862308838
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.742066666666667
$ws.Cells.Item(4,14).Value = 8.2262
$ws.Cells.Item(4,15).Value = 0.1923540942798824
$ws.Cells.Item(4,16).Value = 0.1923540942798823
$ws.Cells.Item(4,17).Value = 2.399414359911111
$ws.Cells.Item(4,18).Value = 21.5947292392
$ws.Cells.Item(4,19).Value = 0.010202574384666
$ws.Cells.Item(4,20).Value = 0.010202574384666

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Ntn1"
$ws.Cells.Item(5,3).Value = "Adora2b"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 10.61942
$ws.Cells.Item(5,8).Value = 31.85826
$ws.Cells.Item(5,9).Value = 0.6436975977089742
$ws.Cells.Item(5,10).Value = 0.6436975977089742
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.888768333333333
$ws.Cells.Item(5,14).Value = 11.666305
$ws.Cells.Item(5,15).Value = 0.2727944290034114
$ws.Cells.Item(5,16).Value = 0.2727944290034114
$ws.Cells.Item(5,17).Value = 41.29646421436667
$ws.Cells.Item(5,18).Value = 371.6681779293
$ws.Cells.Item(5,19).Value = 0.1755971186178872
$ws.Cells.Item(5,20).Value = 0.1755971186178872

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Ntn1"
$ws.Cells.Item(6,3).Value = "Adora2b"
$ws.Cells.Item(6,4).Value = "M2"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 10.61942
$ws.Cells.Item(6,8).Value = 31.85826
$ws.Cells.Item(6,9).Value = 0.6436975977089742
$ws.Cells.Item(6,10).Value = 0.6436975977089742
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 7.624472
$ws.Cells.Item(6,14).Value = 22.873416
$ws.Cells.Item(6,15).Value = 0.5348514767167063
$ws.Cells.Item(6,16).Value = 0.5348514767167063
$ws.Cells.Item(6,17).Value = 80.96747044624
$ws.Cells.Item(6,18).Value = 728.70723401616
$ws.Cells.Item(6,19).Value = 0.3442826106936412
$ws.Cells.Item(6,20).Value = 0.3442826106936412

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Ntn1"
$ws.Cells.Item(7,3).Value = "Adora2b"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 10.61942
$ws.Cells.Item(7,8).Value = 31.85826
$ws.Cells.Item(7,9).Value = 0.6436975977089742
$ws.Cells.Item(7,10).Value = 0.6436975977089742
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.742066666666667
$ws.Cells.Item(7,14).Value = 8.2262
$ws.Cells.Item(7,15).Value = 0.1923540942798824
$ws.Cells.Item(7,16).Value = 0.1923540942798823
$ws.Cells.Item(7,17).Value = 29.11915760133333
$ws.Cells.Item(7,18).Value = 262.072418412
$ws.Cells.Item(7,19).Value = 0.1238178683974458
$ws.Cells.Item(7,20).Value = 0.1238178683974458

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Ntn1"
$ws.Cells.Item(8,3).Value = "Adora2b"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.4029073333333333
$ws.Cells.Item(8,8).Value = 1.208722
$ws.Cells.Item(8,9).Value = 0.02442228319117198
$ws.Cells.Item(8,10).Value = 0.02442228319117198
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 3.888768333333333
$ws.Cells.Item(8,14).Value = 11.666305
$ws.Cells.Item(8,15).Value = 0.2727944290034114
$ws.Cells.Item(8,16).Value = 0.2727944290034114
$ws.Cells.Item(8,17).Value = 1.566813279134444
$ws.Cells.Item(8,18).Value = 14.10131951221
$ws.Cells.Item(8,19).Value = 0.006662262798095371
$ws.Cells.Item(8,20).Value = 0.006662262798095371

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Ntn1"
$ws.Cells.Item(9,3).Value = "Adora2b"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.4029073333333333
$ws.Cells.Item(9,8).Value = 1.208722
$ws.Cells.Item(9,9).Value = 0.02442228319117198
$ws.Cells.Item(9,10).Value = 0.02442228319117198
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 7.624472
$ws.Cells.Item(9,14).Value = 22.873416
$ws.Cells.Item(9,15).Value = 0.5348514767167063
$ws.Cells.Item(9,16).Value = 0.5348514767167063
$ws.Cells.Item(9,17).Value = 3.071955681594666
$ws.Cells.Item(9,18).Value = 27.64760113435199
$ws.Cells.Item(9,19).Value = 0.01306229422959193
$ws.Cells.Item(9,20).Value = 0.01306229422959193

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Ntn1"
$ws.Cells.Item(10,3).Value = "Adora2b"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.4029073333333333
$ws.Cells.Item(10,8).Value = 1.208722
$ws.Cells.Item(10,9).Value = 0.02442228319117198
$ws.Cells.Item(10,10).Value = 0.02442228319117198
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.742066666666667
$ws.Cells.Item(10,14).Value = 8.2262
$ws.Cells.Item(10,15).Value = 0.1923540942798824
$ws.Cells.Item(10,16).Value = 0.1923540942798823
$ws.Cells.Item(10,17).Value = 1.104798768488889
$ws.Cells.Item(10,18).Value = 9.943188916399999
$ws.Cells.Item(10,19).Value = 0.004697726163484681
$ws.Cells.Item(10,20).Value = 0.00469772616348468

# Row 11
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Ntn1"
$ws.Cells.Item(11,3).Value = "Adora2b"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 4.600163333333334
$ws.Cells.Item(11,8).Value = 13.80049
$ws.Cells.Item(11,9).Value = 0.2788395304767656
$ws.Cells.Item(11,10).Value = 0.2788395304767655
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 3.888768333333333
$ws.Cells.Item(11,14).Value = 11.666305
$ws.Cells.Item(11,15).Value = 0.2727944290034114
$ws.Cells.Item(11,16).Value = 0.2727944290034114
$ws.Cells.Item(11,17).Value = 17.88896949882778
$ws.Cells.Item(11,18).Value = 161.00072548945
$ws.Cells.Item(11,19).Value = 0.07606587049998859
$ws.Cells.Item(11,20).Value = 0.07606587049998857

# Row 12
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Ntn1"
$ws.Cells.Item(12,3).Value = "Adora2b"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 4.600163333333334
$ws.Cells.Item(12,8).Value = 13.80049
$ws.Cells.Item(12,9).Value = 0.2788395304767656
$ws.Cells.Item(12,10).Value = 0.2788395304767655
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 7.624472
$ws.Cells.Item(12,14).Value = 22.873416
$ws.Cells.Item(12,15).Value = 0.5348514767167063
$ws.Cells.Item(12,16).Value = 0.5348514767167063
$ws.Cells.Item(12,17).Value = 35.07381653042667
$ws.Cells.Item(12,18).Value = 315.66434877384
$ws.Cells.Item(12,19).Value = 0.1491377346424911
$ws.Cells.Item(12,20).Value = 0.1491377346424911

# Row 13
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Ntn1"
$ws.Cells.Item(13,3).Value = "Adora2b"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 4.600163333333334
$ws.Cells.Item(13,8).Value = 13.80049
$ws.Cells.Item(13,9).Value = 0.2788395304767656
$ws.Cells.Item(13,10).Value = 0.2788395304767655
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 2.742066666666667
$ws.Cells.Item(13,14).Value = 8.2262
$ws.Cells.Item(13,15).Value = 0.1923540942798824
$ws.Cells.Item(13,16).Value = 0.1923540942798823
$ws.Cells.Item(13,17).Value = 12.61395453755556
$ws.Cells.Item(13,18).Value = 113.525590838
$ws.Cells.Item(13,19).Value = 0.0536359253342859
$ws.Cells.Item(13,20).Value = 0.05363592533428588
